# Integrate the mailing list into the webpages: collapse the
# run-per-word Title / Author / Abstract paragraphs into single runs
# holding the full text.

$d = $word.ActiveDocument

$d.Content.Find.Execute("Answers: Completing the square", $false, $false, $false, $false, $false, $true, 1, $false, "Answers: Completing the square", 2)

$d.Content.Find.Execute("Tom Coleman", $false, $false, $false, $false, $false, $true, 1, $false, "Tom Coleman", 2)

$d.Content.Find.Execute("Answers to questions relating to the guide on completing the square.", $false, $false, $false, $false, $false, $true, 1, $false, "Answers to questions relating to the guide on completing the square.", 2)
